# Generate Report for Handoff
# Updates the localization-status report: the 8cfa6c82... file has finished
# translation and is now "Ready for handoff" (machine translation priority),
# with fresh handoff timestamps recorded on the per-locale sheets and
# summarized on the Overview sheet.

$wb = $excel.ActiveWorkbook

# --- Overview sheet: row for 8cfa6c82-ee39-40cd-bdb2-492e1556df60.md ---
$ov = $wb.Worksheets.Item("Overview")
$ov.Range("E3").Value = "Ready for handoff"
$ov.Range("F3").Value = "Ready for handoff"
$ov.Range("G3").Value = "2016-08-17 00:13:11"

# --- zh-cn sheet: row for 8cfa6c82-ee39-40cd-bdb2-492e1556df60.md ---
$zh = $wb.Worksheets.Item("zh-cn")
$zh.Range("C3").Value = "Ready for handoff"
$zh.Range("E3").Value = "mt"
$zh.Range("H3").Value = "2016-08-17 00:13:06"

# --- de-de sheet: row for 8cfa6c82-ee39-40cd-bdb2-492e1556df60.md ---
$de = $wb.Worksheets.Item("de-de")
$de.Range("C3").Value = "Ready for handoff"
$de.Range("E3").Value = "mt"
$de.Range("H3").Value = "2016-08-17 00:13:11"

# --- Column widths widened to fit the longer "Ready for handoff" status text ---
$ov.Columns.Item(5).ColumnWidth = 16.3
$ov.Columns.Item(6).ColumnWidth = 16.3
$zh.Columns.Item(3).ColumnWidth = 16.3
$de.Columns.Item(3).ColumnWidth = 16.3
